$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 471, shifting existing rows 471:564 down to 473:566.
$ws.Range("A471:R472").EntireRow.Insert()

# Populate new row 471 (Calidad = Primera)
$ws.Cells.Item(471, 1).Value = 1
$ws.Cells.Item(471, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(471, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(471, 4).Value = 45275
$ws.Cells.Item(471, 5).Value = 15
$ws.Cells.Item(471, 6).Value = 100114014
$ws.Cells.Item(471, 7).Value = "Betarraga"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 500
$ws.Cells.Item(471, 11).Value = 300
$ws.Cells.Item(471, 12).Value = 400
$ws.Cells.Item(471, 13).Value = 350
$ws.Cells.Item(471, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(471, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(471, 16).Value = 88
$ws.Cells.Item(471, 17).Value = 4
$ws.Cells.Item(471, 18).Value = "Hortaliza"

# Populate new row 472 (Calidad = Segunda)
$ws.Cells.Item(472, 1).Value = 1
$ws.Cells.Item(472, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(472, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(472, 4).Value = 45275
$ws.Cells.Item(472, 5).Value = 15
$ws.Cells.Item(472, 6).Value = 100114014
$ws.Cells.Item(472, 7).Value = "Betarraga"
$ws.Cells.Item(472, 8).Value = "Sin especificar"
$ws.Cells.Item(472, 9).Value = "Segunda"
$ws.Cells.Item(472, 10).Value = 400
$ws.Cells.Item(472, 11).Value = 300
$ws.Cells.Item(472, 12).Value = 400
$ws.Cells.Item(472, 13).Value = 350
$ws.Cells.Item(472, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(472, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(472, 16).Value = 70
$ws.Cells.Item(472, 17).Value = 5
$ws.Cells.Item(472, 18).Value = "Hortaliza"
